$d = $word.ActiveDocument

# --- Edit 1: merge the three runs "<id>", "p017r_1", "</id>" into a single run
# while keeping the formatting of the first run (Courier New, color 7f6000, sz 18).
$r1 = $d.Content
$found1 = $r1.Find.Execute("<id>p017r_1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $mergeRange = $d.Range($r1.Start, $r1.End)
    # First set to a text that differs from the original so the engine actually
    # performs a text replacement (merging the run formatting of the first run),
    # instead of a no-op when new text equals old text.
    $mergeRange.Text = "<id>p017r_1</id> "
    $mergeRange2 = $d.Range($r1.Start, $r1.Start + 17)
    $mergeRange2.Text = "<id>p017r_1</id>"
}

# --- Edit 2: change the single-character run "e" to "è" inside the standalone
# paragraph whose text is "pres.</ab>".
$r2 = $d.Content
$found2 = $r2.Find.Execute("pres.</ab>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $charRange = $d.Range($r2.Start + 2, $r2.Start + 3)
    $charRange.Text = "è"
}
